# Update the three-digit-by-one-digit multiplication problems in the
# single 5-column table. Each problem row in the document corresponds to
# Word table rows 1, 5, 10, 15, 20 (rows in between are blank spacer rows).
# We target each cell directly by (row, column) and overwrite its text,
# which avoids any ambiguity from find/replace collisions (e.g. a newly
# written "828×8=" must NOT be re-matched by a later rule).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Old = "350×4="; New = "254×5=" },
    @{ Row = 1;  Col = 2; Old = "709×5="; New = "628×9=" },
    @{ Row = 1;  Col = 3; Old = "636×6="; New = "806×9=" },
    @{ Row = 1;  Col = 4; Old = "914×7="; New = "991×9=" },
    @{ Row = 1;  Col = 5; Old = "415×4="; New = "200×3=" },

    @{ Row = 5;  Col = 1; Old = "828×8="; New = "942×2=" },
    @{ Row = 5;  Col = 2; Old = "418×6="; New = "704×3=" },
    @{ Row = 5;  Col = 3; Old = "984×7="; New = "347×9=" },
    @{ Row = 5;  Col = 4; Old = "227×7="; New = "438×4=" },
    @{ Row = 5;  Col = 5; Old = "153×6="; New = "422×4=" },

    @{ Row = 10; Col = 1; Old = "379×5="; New = "114×4=" },
    @{ Row = 10; Col = 2; Old = "596×2="; New = "758×2=" },
    @{ Row = 10; Col = 3; Old = "523×6="; New = "518×2=" },
    @{ Row = 10; Col = 4; Old = "240×3="; New = "194×9=" },
    @{ Row = 10; Col = 5; Old = "800×6="; New = "828×8=" },

    @{ Row = 15; Col = 1; Old = "830×8="; New = "260×5=" },
    @{ Row = 15; Col = 2; Old = "164×7="; New = "135×4=" },
    @{ Row = 15; Col = 3; Old = "227×2="; New = "627×9=" },
    @{ Row = 15; Col = 4; Old = "144×5="; New = "644×4=" },
    @{ Row = 15; Col = 5; Old = "279×8="; New = "251×6=" },

    @{ Row = 20; Col = 1; Old = "508×5="; New = "820×7=" },
    @{ Row = 20; Col = 2; Old = "517×5="; New = "609×9=" },
    @{ Row = 20; Col = 3; Old = "592×9="; New = "674×7=" },
    @{ Row = 20; Col = 4; Old = "702×5="; New = "293×7=" },
    @{ Row = 20; Col = 5; Old = "345×6="; New = "163×8=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
